# 13.04.21 debugging (IsoPatternPeakView: copying table, OccupancyRecalculator.py: loading Occupancies_in.csv)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("analysis")

# Update the timestamp label in A1
$ws.Range("A1").Value = "13/04/2021 17:39"

# Updated base-pairing totals
$ws.Range("C3").Value = 0.4967469870553322
$ws.Range("C4").Value = 0.5032530129446677

# Updated occupancy data table (columns C = 5'-O /%, D = 3'-O /%)
$ws.Range("D9").Value  = 0.9702968525780499
$ws.Range("D10").Value = 0.859310439376442

$ws.Range("C11").Value = 0.1724538867617657
$ws.Range("D11").Value = 0.9002876039349853

$ws.Range("D12").Value = 0.952367657234294

$ws.Range("C13").Value = 0.09128218680418146

$ws.Range("C14").Value = 0.08061230026340492
$ws.Range("D14").Value = 0.9089060143233461

$ws.Range("C15").Value = 0.1229420682206724
$ws.Range("D15").Value = 0.8096547457916253
$ws.Range("D15").NumberFormat = "0.0%"

$ws.Range("C16").Value = 0.1235458202654845
$ws.Range("D16").Value = 0.8748118411885112

$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0.90816800837282

$ws.Range("C18").Value = 0.0735704151443306
$ws.Range("D18").Value = 0.8253846146228404

$ws.Range("D19").Value = 1

$ws.Range("C20").Value = 0.5621945415081315
$ws.Range("D20").Value = 0.7748609036594906

$ws.Range("C21").Value = 0.1809350292162644
$ws.Range("D21").Value = 0.8123476569107348

$ws.Range("C22").Value = 0.1636447119069271
$ws.Range("D22").Value = 0.8059151036038091

$ws.Range("C23").Value = 0.1381611054857339
$ws.Range("D23").Value = 0.7574638473677997

$ws.Range("C24").Value = 0.2282825262821861
$ws.Range("D24").Value = 0.8170294457232519

$ws.Range("C25").Value = 0.2276192464347997
$ws.Range("D25").Value = 0.7904972105057587

$ws.Range("C26").Value = 0.3663844860194911
$ws.Range("D26").Value = 0.6868570095767319

$ws.Range("C27").Value = 0.3652077072120486
$ws.Range("D27").Value = 0.678462449211344

$ws.Range("C28").Value = 0.4065624587608148
$ws.Range("D28").Value = 0.7349813382823553

$ws.Range("C29").Value = 0.3750525719654789
$ws.Range("D29").Value = 0.6364863328227983

$ws.Range("C30").Value = 0.402773541720006
$ws.Range("D30").Value = 0.7398683720896356

$ws.Range("C31").Value = 0.4116730285605797
$ws.Range("D31").Value = 0.7322917929956112

$ws.Range("C32").Value = 0.2284720697855734
$ws.Range("D32").Value = 0.6812152183654794

$ws.Range("C33").Value = 0.4324044374920852
$ws.Range("D33").Value = 0.6872980780634957

# Row 34 no longer carries a C value - clear content + formatting entirely
$ws.Range("C34").Clear()
